$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code")

# ---------------------------------------------------------------------------
# 1) Insert a new row at the very top; this shifts every existing row down
#    by one (old row1 -> row2, old row2 -> row3, ... old row9 -> row10).
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new / changed cell values in the same order the original
#    author typed them (this keeps the shared-string table ordering
#    faithful to the source workbook).
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = "component_weight"
$ws.Range("G2").Value = "Python variable"
$ws.Range("F2").Value = "Math Notation"
$ws.Range("H2").Value = "Term"
$ws.Range("I2").Value = "Reason"
$ws.Range("G9").Value = "component_index"
$ws.Range("G5").Value = "spectral_loading"
$ws.Range("D6").Value = "``"
$ws.Range("G6").Value = "positive_sum"
$ws.Range("G7").Value = "negative_sum"
$ws.Range("G3").Value = "data"
$ws.Range("H3").Value = "spectral data"
$ws.Range("F1").Value = "Francis' suggestions"
$ws.Range("I3").Value = "Use what's present in most of the existing manuscript & figures."
$ws.Range("H5").Value = "spectral component loadings"
$ws.Range("H4").Value = "component concentration scores"
$ws.Range("H6").Value = "positive contributions"
$ws.Range("H7").Value = "negative contributions"
$ws.Range("H8").Value = "central contribution offset"
$ws.Range("F9").Value = "i"
$ws.Range("G10").Value = "component_index"

# ---------------------------------------------------------------------------
# 3) Formatting: bold header rows (1 and 2, columns F:I)
# ---------------------------------------------------------------------------
$ws.Range("F1:H1").Font.Bold = $true
$ws.Range("F2:I2").Font.Bold = $true
$ws.Range("F2").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Formatting: wrap text for the new data cells (matches the existing
#    "wrapText" style already used throughout columns A:C)
# ---------------------------------------------------------------------------
$ws.Range("F3:H3").WrapText = $true
$ws.Range("F4:H5").WrapText = $true
$ws.Range("D6").WrapText = $true
$ws.Range("F6").WrapText = $true
$ws.Range("H6:H8").WrapText = $true

# ---------------------------------------------------------------------------
# 5) Merge I3:I5 and center + wrap that merged cell
# ---------------------------------------------------------------------------
$ws.Range("I3:I5").Merge()
$ws.Range("I3:I5").HorizontalAlignment = -4108
$ws.Range("I3:I5").WrapText = $true

# ---------------------------------------------------------------------------
# 6) Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 129.6

# ---------------------------------------------------------------------------
# 7) Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.4987
$ws.Columns.Item(6).ColumnWidth = 13.3854
$ws.Columns.Item(7).ColumnWidth = 18.6081
$ws.Columns.Item(8).ColumnWidth = 19.6081
$ws.Columns.Item(9).ColumnWidth = 22.944

# ---------------------------------------------------------------------------
# 8) Selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("F6").Select()
